$d = $word.ActiveDocument

# --- Edit 1 -----------------------------------------------------------
# "Multiple studies and case reports have implicated ..." -> "Multiple studies implicated ..."
$d.Content.Find.Execute(
    "Multiple studies and case reports have implicated renal artery FMD in impaired kidney function",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Multiple studies implicated renal artery FMD in impaired kidney function", 2) | Out-Null

# --- Edit 2 -----------------------------------------------------------
# Collapse "causal effects of: FMD on chronic kidney disease. We used publicly
# available summary statistics in a two-sample Mendelian randomization study.
# Specifically, we used genetic instruments for FMD ..." into the new wording.
$d.Content.Find.Execute(
    "causal effects of: FMD on chronic kidney disease. We used publicly available summary statistics in a two-sample Mendelian randomization study. Specifically, we used genetic instruments for FMD",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "causal effects of FMD on chronic kidney disease (CKD) in a two-sample Mendelian randomization study with GWAS summary statistics. We used genetic instruments for FMD", 2) | Out-Null

# --- Edit 3 -----------------------------------------------------------
# Remove the three p-value-threshold math equations (and surrounding
# parenthetical text), then reword the trailing sentence.

# Locate the "(" that opens the p-value-threshold parenthetical, and the ")"
# that closes it, and delete everything between (and including) them -- this
# also removes the three embedded <m:oMath> objects. The leading space before
# "(" is consumed too, so only the single space from ") when" remains.
$openRange = $d.Content
$openRange.Find.Execute("We considered three different p-value thresholds (", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$openParenStart = $openRange.End - 2

$closeRange = $d.Content
$closeRange.Find.Execute(") when choosing relevant SNP instruments", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$closeParenEnd = $closeRange.Start + 1

$mathSpan = $d.Range($openParenStart, $closeParenEnd)
$mathSpan.Text = ""

# Reword the remainder of the sentence.
$d.Content.Find.Execute(
    "We considered three different p-value thresholds when choosing relevant SNP instruments for the MR analyses. All five methods showed no evidence of a causal effect of FMD on chronic kidney disease. To further investigate our findings, we performed sensitivity analyses in efforts to assess evidence of horizontal pleiotropy and other sources of confounding",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "We considered three different p-value thresholds when choosing relevant SNP instruments for MR analyses. All five methods showed no evidence of a causal effect of FMD on CKD. To further evaluate our findings, we performed sensitivity analyses to assess evidence of horizontal pleiotropy and other sources of confounding", 2) | Out-Null

Write-Host $d.Paragraphs.Item(5).Range.Text
